# Add season record columns (Wins / Losses / Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column headers in AD1:AF1 ---
$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# Match the header formatting (bold, bordered, centered) used by the
# existing header cells (e.g. A1) by copying their format onto the new
# header cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (rows 2-55): season record for every player row ---
$lastRow = 55
$ws.Range("AD2:AD$lastRow").Value2 = 74
$ws.Range("AE2:AE$lastRow").Value2 = 88
$ws.Range("AF2:AF$lastRow").Value2 = 0
